$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 151.03847
$ws.Range("I2").Value = 156.125
$ws.Range("K2").Value = 156.125
$ws.Range("M2").Value = -43.125
$ws.Range("H40").Value = 5258.2666
$ws.Range("I40").Value = 3374.889
$ws.Range("K40").Value = 3374.889
$ws.Range("M40").Value = -3199.889
$ws.Range("H53").Value = 1420.5333
$ws.Range("J53").Value = 428.625
$ws.Range("L53").Value = 428.625
$ws.Range("N53").Value = -1702.625
$ws.Range("H111").Value = 2150
$ws.Range("I111").Value = 1800
$ws.Range("K111").Value = 5400
$ws.Range("M111").Value = -2333
$ws.Range("H113").Value = 76925970
$ws.Range("I113").Value = 142859730
$ws.Range("K113").Value = 142859730
$ws.Range("M113").Value = -142856476
$ws.Range("H132").Value = 6204.95
$ws.Range("I132").Value = 6473.3687
$ws.Range("K132").Value = 19420.1061
$ws.Range("M132").Value = -16890.1061
$ws.Range("H137").Value = 1429.8572
$ws.Range("I137").Value = 1145.0714
$ws.Range("K137").Value = 3435.2142
$ws.Range("M137").Value = -885.2142000000003
$ws.Range("H138").Value = 1947.2354
$ws.Range("I138").Value = 1135.3103
$ws.Range("J138").Value = 3017.5
$ws.Range("K138").Value = 3405.9309
$ws.Range("L138").Value = 9052.5
$ws.Range("M138").Value = 1734.0691
$ws.Range("N138").Value = -19332.5
$ws.Range("H141").Value = 8592.375
$ws.Range("I141").Value = 7081.9165
$ws.Range("K141").Value = 21245.7495
$ws.Range("M141").Value = -16065.7495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3905.96
$ws.Range("J2").Value = 4925.6
$ws.Range("L2").Value = 4925.6
$ws.Range("N2").Value = -5151.6
$ws.Range("H5").Value = 108.166664
$ws.Range("I5").Value = 108.166664
$ws.Range("K5").Value = 108.166664
$ws.Range("M5").Value = 3.833336000000003
$ws.Range("H32").Value = 4317.2334
$ws.Range("I32").Value = 4417.8447
$ws.Range("K32").Value = 4417.8447
$ws.Range("M32").Value = -4130.8447
$ws.Range("H45").Value = 4155.0225
$ws.Range("I45").Value = 3004
$ws.Range("J45").Value = 6457.067
$ws.Range("K45").Value = 3004
$ws.Range("L45").Value = 6457.067
$ws.Range("M45").Value = -2627
$ws.Range("N45").Value = -7211.067
$ws.Range("H50").Value = 8074
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H61").Value = 3092436.8
$ws.Range("I61").Value = 4171471
$ws.Range("K61").Value = 4171471
$ws.Range("M61").Value = -4171259
$ws.Range("H110").Value = 3571.52
$ws.Range("I110").Value = 2334.0527
$ws.Range("J110").Value = 7490.1665
$ws.Range("K110").Value = 2334.0527
$ws.Range("L110").Value = 7490.1665
$ws.Range("M110").Value = -289.0527000000002
$ws.Range("N110").Value = -11580.1665
$ws.Range("H116").Value = 3905.96
$ws.Range("J116").Value = 4925.6
$ws.Range("L116").Value = 4925.6
$ws.Range("N116").Value = -9513.6
$ws.Range("H122").Value = 3379
$ws.Range("I122").Value = 3379
$ws.Range("K122").Value = 10137
$ws.Range("M122").Value = -7687
$ws.Range("H132").Value = 3463.4133
$ws.Range("I132").Value = 2800.8872
$ws.Range("K132").Value = 8402.661599999999
$ws.Range("M132").Value = -5872.661599999999
$ws.Range("H136").Value = 3092436.8
$ws.Range("I136").Value = 4171471
$ws.Range("K136").Value = 12514413
$ws.Range("M136").Value = -12511863

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3905.96
$ws.Range("J3").Value = 4925.6
$ws.Range("L3").Value = 4925.6
$ws.Range("N3").Value = -5153.6
$ws.Range("H4").Value = 108.166664
$ws.Range("I4").Value = 108.166664
$ws.Range("K4").Value = 108.166664
$ws.Range("M4").Value = 6.833336000000003
$ws.Range("H64").Value = 265.4
$ws.Range("J64").Value = 281.75
$ws.Range("L64").Value = 281.75
$ws.Range("N64").Value = -731.75
$ws.Range("H67").Value = 265.4
$ws.Range("J67").Value = 281.75
$ws.Range("L67").Value = 281.75
$ws.Range("N67").Value = -1841.75
$ws.Range("H105").Value = 4602.385
$ws.Range("I105").Value = 4637.727
$ws.Range("J105").Value = 4408
$ws.Range("K105").Value = 4637.727
$ws.Range("L105").Value = 4408
$ws.Range("M105").Value = -2890.727
$ws.Range("N105").Value = -7902
$ws.Range("H107").Value = 2827.946
$ws.Range("I107").Value = 2338.037
$ws.Range("J107").Value = 4150.7
$ws.Range("K107").Value = 2338.037
$ws.Range("L107").Value = 4150.7
$ws.Range("M107").Value = -418.0369999999998
$ws.Range("N107").Value = -7990.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 16669667
$ws.Range("H31").Value = 7220.9414
$ws.Range("I31").Value = 5827.4165
$ws.Range("K31").Value = 5827.4165
$ws.Range("M31").Value = -5532.4165
$ws.Range("H34").Value = 7220.9414
$ws.Range("I34").Value = 5827.4165
$ws.Range("K34").Value = 5827.4165
$ws.Range("M34").Value = -5625.4165
$ws.Range("H58").Value = 5736.5264
$ws.Range("I58").Value = 1990.5
$ws.Range("J58").Value = 9898.777
$ws.Range("K58").Value = 1990.5
$ws.Range("L58").Value = 9898.777
$ws.Range("M58").Value = -1787.5
$ws.Range("N58").Value = -10304.777
$ws.Range("H99").Value = 7425
$ws.Range("I99").Value = 6010.5
$ws.Range("K99").Value = 6010.5
$ws.Range("M99").Value = -4512.5
$ws.Range("H122").Value = 3737.8076
$ws.Range("J122").Value = 3621.6667
$ws.Range("L122").Value = 10865.0001
$ws.Range("N122").Value = -15765.0001
$ws.Range("H126").Value = 7425
$ws.Range("I126").Value = 6010.5
$ws.Range("K126").Value = 18031.5
$ws.Range("M126").Value = -15561.5
$ws.Range("H136").Value = 5736.5264
$ws.Range("I136").Value = 1990.5
$ws.Range("J136").Value = 9898.777
$ws.Range("K136").Value = 5971.5
$ws.Range("L136").Value = 29696.331
$ws.Range("M136").Value = -3421.5
$ws.Range("N136").Value = -34796.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1936.1428
$ws.Range("I5").Value = 316.25
$ws.Range("J5").Value = 4096
$ws.Range("K5").Value = 948.75
$ws.Range("L5").Value = 12288
$ws.Range("M5").Value = -836.75
$ws.Range("N5").Value = -12512
$ws.Range("H135").Value = 1936.1428
$ws.Range("I135").Value = 316.25
$ws.Range("J135").Value = 4096
$ws.Range("K135").Value = 2846.25
$ws.Range("L135").Value = 36864
$ws.Range("M135").Value = -311.25
$ws.Range("N135").Value = -41934

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -48494
$ws.Range("H126").Value = 3030.0476
$ws.Range("J126").Value = 3630.4
$ws.Range("L126").Value = 10891.2
$ws.Range("N126").Value = -15831.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1524.875
$ws.Range("I13").Value = 2399.6667
$ws.Range("K13").Value = 2399.6667
$ws.Range("M13").Value = -2259.6667
$ws.Range("H22").Value = 1162.5264
$ws.Range("I22").Value = 848.5
$ws.Range("J22").Value = 1511.4445
$ws.Range("K22").Value = 848.5
$ws.Range("L22").Value = 1511.4445
$ws.Range("M22").Value = -553.5
$ws.Range("N22").Value = -2101.4445
$ws.Range("H27").Value = 1162.5264
$ws.Range("I27").Value = 848.5
$ws.Range("J27").Value = 1511.4445
$ws.Range("K27").Value = 848.5
$ws.Range("L27").Value = 1511.4445
$ws.Range("M27").Value = -741.5
$ws.Range("N27").Value = -1725.4445
$ws.Range("H132").Value = 13106.143
$ws.Range("I132").Value = 16864.3
$ws.Range("J132").Value = 3710.75
$ws.Range("K132").Value = 50592.89999999999
$ws.Range("L132").Value = 11132.25
$ws.Range("M132").Value = -48062.89999999999
$ws.Range("N132").Value = -16192.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 2790.9375
$ws.Range("J122").Value = 4452.778
$ws.Range("L122").Value = 13358.334
$ws.Range("N122").Value = -18258.334
$ws.Range("H124").Value = 25000
$ws.Range("J124").Value = 25000
$ws.Range("L124").Value = 25000
$ws.Range("N124").Value = -34820
$ws.Range("H132").Value = 6411.0415
$ws.Range("I132").Value = 5404.125
$ws.Range("K132").Value = 16212.375
$ws.Range("M132").Value = -13682.375
